# Regenerate save_data to use K instead of Strike#, recalc std/mean and write s_vals.
# The "K" column (column G) values are recomputed; write the new values in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 1
    6  = 4
    7  = 3
    8  = 0
    9  = 4
    10 = 0
    11 = 3
    12 = 0
    13 = 1
    14 = 0
    15 = 2
    16 = 2
    17 = 2
    18 = 1
    19 = 2
    20 = 2
    21 = 4
    22 = 4
    23 = 5
    24 = 2
    25 = 0
    26 = 2
    27 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
